# Development-Tools.pptx edit
# 1. Merge the two "Project " / "tools" runs on slide 2 into a single run "Project tools".
# 2. Update the Notes Master's fixed date placeholder text from "06-Oct-14" to "1/3/2016".
# 3. Remove the two trailing "Homework" / "Homework (2)" slides (slide 29 & slide 30).

$p = $ppt.ActivePresentation

# --- 1. Merge "Project " + "tools" runs on slide 2 (Table of Contents) ---
$s2 = $p.Slides.Item(2)
$contentShape = $s2.Shapes.Item(2)
$tr = $contentShape.TextFrame.TextRange
$firstPara = $tr.Paragraphs(1, 1)
# Force a real text delta (avoid a same-text no-op) before writing the final value.
$firstPara.Text = "ZZZ__placeholder__ZZZ"
$firstPara2 = $tr.Paragraphs(1, 1)
$firstPara2.Text = "Project tools"

# --- 2. Notes Master fixed date text ---
$nm = $p.NotesMaster
$nmHF = $nm.HeadersFooters
$nmHF.DateAndTime.UseFormat = 0
$nmHF.DateAndTime.Text = "1/3/2016"

# --- 3. Delete the trailing "Homework" / "Homework (2)" slides ---
for ($i = $p.Slides.Count; $i -ge 1; $i--) {
    $candidate = $p.Slides.Item($i)
    $candidateTitle = $candidate.Shapes.Item(1).TextFrame.TextRange.Text
    if ($candidateTitle -like "Homework*") {
        $candidate.Delete()
    }
}
